# Remove the "XOP" / "commodities" holding row from the Holdings sheet.
# This shifts all subsequent rows up by one (row 34->33, 35->34, ... 44->43),
# shrinking the used range from A1:B44 to A1:B43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holdings")

# Row 33 currently holds "XOP" (Index) / "commodities" (Set) - delete it entirely.
$ws.Rows.Item(33).Delete()
